$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 16550
$ws.Cells.Item(34, 10).Value = 25500
$ws.Cells.Item(34, 12).Value = 25500
$ws.Cells.Item(34, 14).Value = -25906
$ws.Cells.Item(36, 8).Value = 16550
$ws.Cells.Item(36, 10).Value = 25500
$ws.Cells.Item(36, 12).Value = 25500
$ws.Cells.Item(36, 14).Value = -26930
$ws.Cells.Item(61, 8).Value = 574.25
$ws.Cells.Item(61, 9).Value = 93.333336
$ws.Cells.Item(61, 10).Value = 2017
$ws.Cells.Item(61, 11).Value = 280.000008
$ws.Cells.Item(61, 12).Value = 6051
$ws.Cells.Item(61, 13).Value = -108.000008
$ws.Cells.Item(61, 14).Value = -6395
$ws.Cells.Item(88, 8).Value = 27780378
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 27780378
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 27780378
$ws.Cells.Item(88, 13).Value = $null
$ws.Cells.Item(88, 14).Value = -27781190
$ws.Cells.Item(91, 8).Value = 27780378
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 27780378
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 27780378
$ws.Cells.Item(91, 13).Value = $null
$ws.Cells.Item(91, 14).Value = -27783186
$ws.Cells.Item(98, 8).Value = 1242946.4
$ws.Cells.Item(98, 9).Value = 1597631
$ws.Cells.Item(98, 10).Value = 1550
$ws.Cells.Item(98, 11).Value = 1597631
$ws.Cells.Item(98, 12).Value = 1550
$ws.Cells.Item(98, 13).Value = -1596133
$ws.Cells.Item(98, 14).Value = -4546
$ws.Cells.Item(122, 8).Value = 1242946.4
$ws.Cells.Item(122, 9).Value = 1597631
$ws.Cells.Item(122, 10).Value = 1550
$ws.Cells.Item(122, 11).Value = 4792893
$ws.Cells.Item(122, 12).Value = 4650
$ws.Cells.Item(122, 13).Value = -4790443
$ws.Cells.Item(122, 14).Value = -9550

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 50000.5
$ws.Cells.Item(76, 10).Value = 50000.5
$ws.Cells.Item(76, 12).Value = 50000.5
$ws.Cells.Item(76, 14).Value = -50676.5
$ws.Cells.Item(79, 8).Value = 50000.5
$ws.Cells.Item(79, 10).Value = 50000.5
$ws.Cells.Item(79, 12).Value = 50000.5
$ws.Cells.Item(79, 14).Value = -52340.5
$ws.Cells.Item(96, 8).Value = 29999.75
$ws.Cells.Item(96, 10).Value = 29999.75
$ws.Cells.Item(96, 12).Value = 29999.75
$ws.Cells.Item(96, 14).Value = -35491.75
$ws.Cells.Item(122, 8).Value = 21891.8
$ws.Cells.Item(122, 9).Value = 34493
$ws.Cells.Item(122, 11).Value = 103479
$ws.Cells.Item(122, 13).Value = -101029

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3557.2415
$ws.Cells.Item(134, 9).Value = 2400.2
$ws.Cells.Item(134, 10).Value = 6128.4443
$ws.Cells.Item(134, 11).Value = 7200.599999999999
$ws.Cells.Item(134, 12).Value = 18385.3329
$ws.Cells.Item(134, 13).Value = -4665.599999999999
$ws.Cells.Item(134, 14).Value = -23455.3329

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 813.55554
$ws.Cells.Item(16, 9).Value = 801.8333
$ws.Cells.Item(16, 10).Value = 837
$ws.Cells.Item(16, 11).Value = 801.8333
$ws.Cells.Item(16, 12).Value = 837
$ws.Cells.Item(16, 13).Value = -514.8333
$ws.Cells.Item(16, 14).Value = -1411
$ws.Cells.Item(113, 8).Value = 813.55554
$ws.Cells.Item(113, 9).Value = 801.8333
$ws.Cells.Item(113, 10).Value = 837
$ws.Cells.Item(113, 11).Value = 801.8333
$ws.Cells.Item(113, 12).Value = 837
$ws.Cells.Item(113, 13).Value = 1368.1667
$ws.Cells.Item(113, 14).Value = -5177
$ws.Cells.Item(122, 8).Value = 1017.8333
$ws.Cells.Item(122, 9).Value = 861.4
$ws.Cells.Item(122, 10).Value = 1800
$ws.Cells.Item(122, 11).Value = 2584.2
$ws.Cells.Item(122, 12).Value = 5400
$ws.Cells.Item(122, 13).Value = -134.1999999999998
$ws.Cells.Item(122, 14).Value = -10300

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 258.4
$ws.Cells.Item(98, 9).Value = 346
$ws.Cells.Item(98, 10).Value = 200
$ws.Cells.Item(98, 11).Value = 1038
$ws.Cells.Item(98, 12).Value = 600
$ws.Cells.Item(98, 13).Value = 460
$ws.Cells.Item(98, 14).Value = -3596

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 9389556
$ws.Cells.Item(11, 9).Value = 12071429
$ws.Cells.Item(11, 10).Value = 3000
$ws.Cells.Item(11, 11).Value = 12071429
$ws.Cells.Item(11, 12).Value = 3000
$ws.Cells.Item(11, 13).Value = -12071290
$ws.Cells.Item(11, 14).Value = -3278
$ws.Cells.Item(12, 8).Value = 3500
$ws.Cells.Item(12, 9).Value = 2000
$ws.Cells.Item(12, 10).Value = 5000
$ws.Cells.Item(12, 11).Value = 2000
$ws.Cells.Item(12, 12).Value = 5000
$ws.Cells.Item(12, 13).Value = -1860
$ws.Cells.Item(12, 14).Value = -5280
$ws.Cells.Item(18, 8).Value = 7500
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 7500
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 7500
$ws.Cells.Item(18, 13).Value = $null
$ws.Cells.Item(18, 14).Value = -8086
$ws.Cells.Item(70, 8).Value = 5545.3335
$ws.Cells.Item(70, 9).Value = 5995.8
$ws.Cells.Item(70, 10).Value = 4258.2856
$ws.Cells.Item(70, 11).Value = 5995.8
$ws.Cells.Item(70, 12).Value = 4258.2856
$ws.Cells.Item(70, 13).Value = -5725.8
$ws.Cells.Item(70, 14).Value = -4798.2856
$ws.Cells.Item(73, 8).Value = 5545.3335
$ws.Cells.Item(73, 9).Value = 5995.8
$ws.Cells.Item(73, 10).Value = 4258.2856
$ws.Cells.Item(73, 11).Value = 5995.8
$ws.Cells.Item(73, 12).Value = 4258.2856
$ws.Cells.Item(73, 13).Value = -5059.8
$ws.Cells.Item(73, 14).Value = -6130.2856
$ws.Cells.Item(107, 8).Value = 987
$ws.Cells.Item(107, 9).Value = 1447.4286
$ws.Cells.Item(107, 10).Value = 449.83334
$ws.Cells.Item(107, 11).Value = 1447.4286
$ws.Cells.Item(107, 12).Value = 449.83334
$ws.Cells.Item(107, 13).Value = 472.5714
$ws.Cells.Item(107, 14).Value = -4289.83334
$ws.Cells.Item(122, 8).Value = 1236890.1
$ws.Cells.Item(122, 9).Value = 1853835.1
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 5561505.300000001
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -5559055.300000001
$ws.Cells.Item(122, 14).Value = -13900

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2974.5715
$ws.Cells.Item(7, 9).Value = 1982.2727
$ws.Cells.Item(7, 11).Value = 1982.2727
$ws.Cells.Item(7, 13).Value = -1870.2727
$ws.Cells.Item(13, 8).Value = 20000
$ws.Cells.Item(25, 8).Value = 50000
$ws.Cells.Item(25, 9).Value = 50000
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 50000
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = -49770
$ws.Cells.Item(25, 14).Value = $null
$ws.Cells.Item(40, 8).Value = 3974.95
$ws.Cells.Item(40, 9).Value = 1933
$ws.Cells.Item(40, 10).Value = 4335.294
$ws.Cells.Item(40, 11).Value = 1933
$ws.Cells.Item(40, 12).Value = 4335.294
$ws.Cells.Item(40, 13).Value = -1797
$ws.Cells.Item(40, 14).Value = -4607.294
$ws.Cells.Item(46, 8).Value = 2174.4
$ws.Cells.Item(46, 9).Value = 1900
$ws.Cells.Item(46, 10).Value = 2243
$ws.Cells.Item(46, 11).Value = 1900
$ws.Cells.Item(46, 12).Value = 2243
$ws.Cells.Item(46, 13).Value = -1712
$ws.Cells.Item(46, 14).Value = -2619
$ws.Cells.Item(122, 8).Value = 3832.9524
$ws.Cells.Item(122, 9).Value = 3273
$ws.Cells.Item(122, 10).Value = 3964.7058
$ws.Cells.Item(122, 11).Value = 9819
$ws.Cells.Item(122, 12).Value = 11894.1174
$ws.Cells.Item(122, 13).Value = -7369
$ws.Cells.Item(122, 14).Value = -16794.1174
$ws.Cells.Item(126, 8).Value = 2974.5715
$ws.Cells.Item(126, 9).Value = 1982.2727
$ws.Cells.Item(126, 11).Value = 5946.8181
$ws.Cells.Item(126, 13).Value = -3476.8181
$ws.Cells.Item(132, 8).Value = 4409.4443
$ws.Cells.Item(132, 9).Value = 3284.1
$ws.Cells.Item(132, 11).Value = 9852.299999999999
$ws.Cells.Item(132, 13).Value = -7322.299999999999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 335333.34
$ws.Cells.Item(81, 9).Value = 500500
$ws.Cells.Item(81, 10).Value = 5000
$ws.Cells.Item(81, 11).Value = 1001000
$ws.Cells.Item(81, 12).Value = 10000
$ws.Cells.Item(81, 13).Value = -999939
$ws.Cells.Item(81, 14).Value = -12122
$ws.Cells.Item(84, 8).Value = 335333.34
$ws.Cells.Item(84, 9).Value = 500500
$ws.Cells.Item(84, 10).Value = 5000
$ws.Cells.Item(84, 11).Value = 5005000
$ws.Cells.Item(84, 12).Value = 50000
$ws.Cells.Item(84, 13).Value = -4999696
$ws.Cells.Item(84, 14).Value = -60608
$ws.Cells.Item(122, 8).Value = 64393.938
$ws.Cells.Item(122, 9).Value = 101280.8
$ws.Cells.Item(122, 10).Value = 2915.8333
$ws.Cells.Item(122, 11).Value = 303842.4
$ws.Cells.Item(122, 12).Value = 8747.499899999999
$ws.Cells.Item(122, 13).Value = -301392.4
$ws.Cells.Item(122, 14).Value = -13647.4999
$ws.Cells.Item(123, 8).Value = 30936.143
$ws.Cells.Item(123, 10).Value = 30936.143
$ws.Cells.Item(123, 12).Value = 30936.143
$ws.Cells.Item(123, 14).Value = -40736.143
